$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 33616.684
$ws.Range("J51").Value = 72577.10000000001
$ws.Range("L51").Value = 72577.10000000001
$ws.Range("N51").Value = -73545.10000000001

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 217.22728
$ws.Range("I55").Value = 306.2857
$ws.Range("J55").Value = 175.66667
$ws.Range("K55").Value = 306.2857
$ws.Range("L55").Value = 175.66667
$ws.Range("M55").Value = -92.28570000000002
$ws.Range("N55").Value = -603.6666700000001

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2521.6316
$ws.Range("I141").Value = 2466.3635
$ws.Range("J141").Value = 2886.4
$ws.Range("K141").Value = 7399.0905
$ws.Range("L141").Value = 8659.200000000001
$ws.Range("M141").Value = -2219.0905
$ws.Range("N141").Value = -19019.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2897.7144
$ws.Range("I61").Value = 2507.8147
$ws.Range("J61").Value = 3599.5334
$ws.Range("K61").Value = 2507.8147
$ws.Range("L61").Value = 3599.5334
$ws.Range("M61").Value = -2295.8147
$ws.Range("N61").Value = -4023.5334

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2897.7144
$ws.Range("I136").Value = 2507.8147
$ws.Range("J136").Value = 3599.5334
$ws.Range("K136").Value = 7523.4441
$ws.Range("L136").Value = 10798.6002
$ws.Range("M136").Value = -4973.4441
$ws.Range("N136").Value = -15898.6002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3042.5
$ws.Range("I94").Value = 1649.5
$ws.Range("K94").Value = 1649.5
$ws.Range("M94").Value = -1198.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20836020
$ws.Range("I134").Value = 2759.111
$ws.Range("K134").Value = 8277.332999999999
$ws.Range("M134").Value = -5742.332999999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 24503.834
$ws.Range("I6").Value = 2673
$ws.Range("J6").Value = 46334.668
$ws.Range("K6").Value = 2673
$ws.Range("L6").Value = 46334.668
$ws.Range("M6").Value = -2560
$ws.Range("N6").Value = -46560.668

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 7143207
$ws.Range("I19").Value = 7692488
$ws.Range("K19").Value = 7692488
$ws.Range("M19").Value = -7692318

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 7143207
$ws.Range("I24").Value = 7692488
$ws.Range("K24").Value = 7692488
$ws.Range("M24").Value = -7692318

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 10047.5
$ws.Range("I36").Value = 10047.5
$ws.Range("K36").Value = 10047.5
$ws.Range("M36").Value = -9659.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 10047.5
$ws.Range("I40").Value = 10047.5
$ws.Range("K40").Value = 10047.5
$ws.Range("M40").Value = -9887.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 49999.5
$ws.Range("J70").Value = 49999.5
$ws.Range("L70").Value = 49999.5
$ws.Range("N70").Value = -50629.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 49999.5
$ws.Range("J73").Value = 49999.5
$ws.Range("L73").Value = 49999.5
$ws.Range("N73").Value = -52183.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1516.3636
$ws.Range("J5").Value = 1404.5
$ws.Range("L5").Value = 4213.5
$ws.Range("N5").Value = -4437.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2977.375
$ws.Range("I69").Value = 1674.6666
$ws.Range("K69").Value = 5023.9998
$ws.Range("M69").Value = -4212.9998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2977.375
$ws.Range("I72").Value = 1674.6666
$ws.Range("K72").Value = 15071.9994
$ws.Range("M72").Value = -11015.9994

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4716.5
$ws.Range("I80").Value = 4399
$ws.Range("K80").Value = 13197
$ws.Range("M80").Value = -12261

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4716.5
$ws.Range("I83").Value = 4399
$ws.Range("K83").Value = 39591
$ws.Range("M83").Value = -34911

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 496.5
$ws.Range("I98").Value = 378.4
$ws.Range("J98").Value = 693.3333
$ws.Range("K98").Value = 1135.2
$ws.Range("L98").Value = 2079.9999
$ws.Range("M98").Value = 362.8000000000002
$ws.Range("N98").Value = -5075.9999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 5350.2
$ws.Range("I134").Value = 5350.2
$ws.Range("K134").Value = 16050.6
$ws.Range("M134").Value = -10980.6

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1516.3636
$ws.Range("J135").Value = 1404.5
$ws.Range("L135").Value = 12640.5
$ws.Range("N135").Value = -17710.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8580971
$ws.Range("I3").Value = 3334600
$ws.Range("J3").Value = 12515748
$ws.Range("K3").Value = 3334600
$ws.Range("L3").Value = 12515748
$ws.Range("M3").Value = -3334484
$ws.Range("N3").Value = -12515980

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 5010720.5
$ws.Range("I14").Value = 6263038
$ws.Range("K14").Value = 6263038
$ws.Range("M14").Value = -6262870

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2722.8333
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1461.2
$ws.Range("I22").Value = 918.5
$ws.Range("J22").Value = 1962.1538
$ws.Range("K22").Value = 918.5
$ws.Range("L22").Value = 1962.1538
$ws.Range("M22").Value = -623.5
$ws.Range("N22").Value = -2552.1538

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1461.2
$ws.Range("I27").Value = 918.5
$ws.Range("J27").Value = 1962.1538
$ws.Range("K27").Value = 918.5
$ws.Range("L27").Value = 1962.1538
$ws.Range("M27").Value = -811.5
$ws.Range("N27").Value = -2176.1538

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2299.342
$ws.Range("J46").Value = 2402.7646
$ws.Range("L46").Value = 2402.7646
$ws.Range("N46").Value = -2778.7646

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 45458040
$ws.Range("I136").Value = 3228.8462
$ws.Range("J136").Value = 111114984
$ws.Range("K136").Value = 9686.5386
$ws.Range("L136").Value = 333344952
$ws.Range("M136").Value = -7136.5386
$ws.Range("N136").Value = -333350052

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 18355.572
$ws.Range("J9").Value = 1100
$ws.Range("L9").Value = 1100
$ws.Range("N9").Value = -1380

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 37699.9
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 37699.9
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 37699.9
$ws.Range("N70").Value = -38329.9
$ws.Range("M70").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 37699.9
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 37699.9
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 37699.9
$ws.Range("N73").Value = -39883.9
$ws.Range("M73").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3231.5
$ws.Range("I122").Value = 3271
$ws.Range("K122").Value = 9813
$ws.Range("M122").Value = -7363
